$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$changes = @(
  @{ Row = 1;  Col = 1; Old = "31÷6="; New = "31÷2=" },
  @{ Row = 1;  Col = 2; Old = "72÷6="; New = "58÷8=" },
  @{ Row = 1;  Col = 3; Old = "88÷7="; New = "44÷7=" },
  @{ Row = 1;  Col = 4; Old = "30÷8="; New = "43÷8=" },
  @{ Row = 1;  Col = 5; Old = "25÷3="; New = "90÷7=" },

  @{ Row = 5;  Col = 1; Old = "69÷2="; New = "20÷7=" },
  @{ Row = 5;  Col = 2; Old = "78÷7="; New = "55÷4=" },
  @{ Row = 5;  Col = 3; Old = "94÷8="; New = "44÷2=" },
  @{ Row = 5;  Col = 4; Old = "56÷3="; New = "72÷6=" },
  @{ Row = 5;  Col = 5; Old = "78÷4="; New = "87÷5=" },

  @{ Row = 9;  Col = 1; Old = "88÷7="; New = "86÷7=" },
  @{ Row = 9;  Col = 2; Old = "13÷8="; New = "63÷9=" },
  @{ Row = 9;  Col = 3; Old = "95÷6="; New = "52÷6=" },
  @{ Row = 9;  Col = 4; Old = "98÷2="; New = "52÷2=" },
  @{ Row = 9;  Col = 5; Old = "59÷5="; New = "36÷4=" },

  @{ Row = 13; Col = 1; Old = "85÷2="; New = "81÷7=" },
  @{ Row = 13; Col = 2; Old = "96÷3="; New = "62÷2=" },
  @{ Row = 13; Col = 3; Old = "87÷9="; New = "97÷9=" },
  @{ Row = 13; Col = 4; Old = "79÷7="; New = "11÷6=" },
  @{ Row = 13; Col = 5; Old = "96÷8="; New = "20÷6=" },

  @{ Row = 17; Col = 1; Old = "20÷4="; New = "20÷2=" },
  @{ Row = 17; Col = 2; Old = "39÷2="; New = "39÷9=" },
  @{ Row = 17; Col = 3; Old = "31÷8="; New = "90÷6=" },
  @{ Row = 17; Col = 4; Old = "94÷9="; New = "49÷3=" },
  @{ Row = 17; Col = 5; Old = "90÷2="; New = "88÷3=" }
)

foreach ($change in $changes) {
  $cell = $t.Cell($change.Row, $change.Col)
  $range = $cell.Range
  # Replace = 1 (wdReplaceOne) so the edit stays confined to this cell's
  # range, since wdReplaceAll (2) here replaces every matching occurrence
  # across the whole document instead of just within the given range.
  [void]$range.Find.Execute($change.Old, $true, $false, $false, $false, $false, $true, 0, $false, $change.New, 1)
}
